# Daily attendance processing - 2026-01-20 23:02:02
# Normalize the "Recorded By" (column G) value ordering on the
# "Session Analysis Results" sheet: for a known set of exact string
# values, rotate the comma-separated list of recorders left by one
# position (first entry moves to the end).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    "system, System, backup@backdoor.com" = "System, backup@backdoor.com, system"
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System"
    "backup@backdoor.com, System"         = "System, backup@backdoor.com"
    "admin@admin.com, dnasr281@gmail.com" = "dnasr281@gmail.com, admin@admin.com"
}

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $current = $cell.Text

    if ($map.ContainsKey($current)) {
        $cell.Value = $map[$current]
    }
}
